# Insert a new data row for the Chirimoya price list: a new week's record
# is inserted at row 81, pushing the existing rows 81-115 down to 82-116.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 81, shifting rows 81:115 down to 82:116.
$ws.Rows("81:81").Insert()

# Populate the new row 81 with the new record's values.
$ws.Cells.Item(81, 1).Value = 10
$ws.Cells.Item(81, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(81, 3).Value = "La Araucanía"
$ws.Cells.Item(81, 4).Value = 44529
$ws.Cells.Item(81, 5).Value = 9
$ws.Cells.Item(81, 6).Value = "Fruta"
$ws.Cells.Item(81, 7).Value = 100107
$ws.Cells.Item(81, 8).Value = "Otros"
$ws.Cells.Item(81, 9).Value = 100107002
$ws.Cells.Item(81, 10).Value = "Chirimoya"
$ws.Cells.Item(81, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(81, 12).Value = "Primera"
$ws.Cells.Item(81, 13).Value = 65
$ws.Cells.Item(81, 14).Value = 3000
$ws.Cells.Item(81, 15).Value = 3000
$ws.Cells.Item(81, 16).Value = 3000
$ws.Cells.Item(81, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(81, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(81, 19).Value = 3000
$ws.Cells.Item(81, 20).Value = 1
